# "9th Stab - Cosmetic Changes"
#
# The sheet is a MarketBeat rank watch: column A = analyst/firm name,
# then one column per observation week (newest week leftmost, right next
# to column A). Before this edit the weeks were "Jun_13" (col B) and
# "Jun_10" (col C). This edit records two more - newer - weekly
# snapshots ("Jun_17", "Jun_15"), pushing the existing weeks to the
# right, so the header becomes:
#   B1=Jun_17 (new)  C1=Jun_15 (new)  D1=Jun_13 (was B1)  E1=Jun_10 (was C1)
# and every data row gets two new "UN" cells in the newly inserted C/D
# columns, with whatever was in the old column C (plain "UN", or the
# special dated/highlighted rating-change note on rows 10 and 17)
# simply sliding two columns right into E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before the existing column C. This shifts the
# old column C (values + styles, e.g. the highlighted note in C17) to
# column E untouched, and leaves C:D blank/new. Column B (and its value)
# is left alone by the insert itself.
$ws.Columns("C:D").Insert()

# New columns C & D get the same "UN" filler the rest of column B/old-C
# already used, for every data row.
$ws.Range("C2:D27").Value = "UN"

# Row 1 header: old B1 ("Jun_13") slides right into the new D1 ...
$ws.Range("D1").Value = $ws.Range("B1").Value2
# ... and B1/C1 become the two brand new (more recent) week labels.
# C1 is written first so the new shared strings land in the same
# left-to-right order as the columns: Jun_15 before Jun_17.
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Match column C's pre-existing custom width (8 chars) on the two new
# columns plus the shifted-to column E.
$ws.Columns("C:E").ColumnWidth = 7.15

# Cosmetic outline/group markers the original diff shows on the new
# column pair (harmless no-op on engines that don't serialize it).
$ws.Columns("C:D").Collapsed = $true
$ws.Columns("E:E").Collapsed = $false
